$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the summation formulas that were added in F7, D8, F8, D9
$ws.Range("F7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D9").ClearContents()

# Update the selected cell on the sheet view
$ws.Range("C17").Select()
